$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new day of Argent (silver) price data as row 72, mirroring the
# existing rows where every column -- including the numeric-looking ones --
# is stored as plain text. Temporarily format the range as Text so Excel
# doesn't auto-convert the date/number-looking strings, then restore the
# Normal style so no stray per-cell formatting is left behind.
$rng = $ws.Range("A72:J72")
$rng.NumberFormat = "@"

$ws.Range("A72").Value = "2025-05-12"
$ws.Range("B72").Value = "38"
$ws.Range("C72").Value = "37.05"
$ws.Range("D72").Value = "0.98"
$ws.Range("E72").Value = "0.265"
$ws.Range("F72").Value = "0.09"
$ws.Range("G72").Value = "5,311"
$ws.Range("H72").Value = "7,951"
$ws.Range("I72").Value = "8,001"
$ws.Range("J72").Value = "7.2617"

$rng.Style = "Normal"
